$wb = $excel.ActiveWorkbook

# 1. Rename the second sheet ("Include from Onkologie Residu" -> "Include #0")
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "Include #0"

# 2. On the Metadata sheet, update the Date value and insert a new
#    "Jurisdiction" property row right after "Contact".
$ws1 = $wb.Worksheets.Item(1)

# Update Date value (row 8, column B)
$ws1.Range("B8").Value = "2024-09-17T19:55:11+00:00"

# Insert a new row after "Contact" (row 10), before "Description" (row 11),
# copying the formatting of an existing body row so styles match.
$ws1.Rows.Item(11).Insert()
$ws1.Range("A10:B10").Copy()
$ws1.Range("A11:B11").PasteSpecial(-4122)

$ws1.Range("A11").Value = "Jurisdiction"
$ws1.Range("B11").Value = ""
